$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B28").Value = 181
$ws.Range("C28").Value = 181
$ws.Range("D28").Value = 181
$ws.Range("E28").Value = 232

[void]$ws.Range("C42").Select()

$wb.Save()
